# Implement data parsing logic
# Appends one new parsed-record row to each of the four log sheets, matching
# newly-decoded entries (time / length / id / actual-length / checksum, plus
# their decimal companions).

$wb = $excel.ActiveWorkbook

function Add-LogRow($ws, $row, $time, $b, $c, $d, $e, $f, $g, $h, $i) {
    $ws.Cells.Item($row, 1).Value = $time
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

$bigId = [double]"5.68631262647114e+23"

# --- ROW50-FE-LIFTER: new row 36 -------------------------------------------
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-LogRow $ws1 36 45742.1642887963 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x72" "0xe" 400 $bigId 370 14

# --- ROW50-MID-LIFTER: new row 38 -------------------------------------------
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-LogRow $ws2 38 45742.13329861111 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x76" "0x19" 400 0 374 25
# G38 in the source is recorded as literal text (the integer exceeds safe
# double precision), so force a text format before writing and then strip the
# explicit format again so the cell keeps the default style.
$ws2.Cells.Item(38, 7).NumberFormat = "@"
$ws2.Cells.Item(38, 7).Value = "568631262647113771663628"
$ws2.Cells.Item(38, 7).ClearFormats()
# ClearFormats also clears the date format on A38, so restore it afterwards.
$ws2.Cells.Item(38, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- ROW11-FE-LIFTER: new row 36 --------------------------------------------
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-LogRow $ws3 36 45742.18000766203 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x72" "0x14" 400 $bigId 370 20

# --- ROW11-MID-LIFTER: new row 36 -------------------------------------------
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-LogRow $ws4 36 45742.3298003125 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x7a" "0x19" 400 $bigId 378 25
